$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 710, shifting existing rows 710:769 down to 711:770
$ws.Rows.Item(710).Insert()

# Populate the new row 710 with the updated record
$ws.Cells.Item(710, 1).Value = 10
$ws.Cells.Item(710, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(710, 3).Value = "La Araucanía"
$ws.Cells.Item(710, 4).Value = 44826
$ws.Cells.Item(710, 5).Value = 9
$ws.Cells.Item(710, 6).Value = 100112003
$ws.Cells.Item(710, 7).Value = "Ajo"
$ws.Cells.Item(710, 8).Value = "Chino"
$ws.Cells.Item(710, 9).Value = "Primera"
$ws.Cells.Item(710, 10).Value = 1050
$ws.Cells.Item(710, 11).Value = 2000
$ws.Cells.Item(710, 12).Value = 24000
$ws.Cells.Item(710, 13).Value = 15762
$ws.Cells.Item(710, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(710, 15).Value = "China"
$ws.Cells.Item(710, 16).Value = 1576
$ws.Cells.Item(710, 17).Value = 10
$ws.Cells.Item(710, 18).Value = "Hortaliza"
